$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the two GiantSpiderRoom occurrences with EmptyCavePath
$ws.Range("D4").Value = "EmptyCavePath"
$ws.Range("B7").Value = "EmptyCavePath"

# Add the new header/instruction row
$ws.Range("A1").Value = "Syntax: tile_name, other attributes"
$ws.Range("A1").HorizontalAlignment = -4131
$ws.Range("A1").VerticalAlignment = -4108
$ws.Rows.Item(1).RowHeight = 45

$ws.Range("C1").Value = "Other attributes are added with syntax `"<object name>:<quantity>`""
$ws.Range("C1").HorizontalAlignment = -4131
$ws.Range("C1").VerticalAlignment = -4160

# Update selection to match the new active cell / selected range
$ws.Range("B3:H8").Select()
